$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Row 21: mark task as done ("Hecho"), total = 1, week column W (6th week pair) = 1
$ws.Range("F21").Value = "Hecho"
$ws.Range("G21").Value = 1
$ws.Range("W21").Value = 1

# Row 22: same as row 21
$ws.Range("F22").Value = "Hecho"
$ws.Range("G22").Value = 1
$ws.Range("W22").Value = 1

# Row 23: total = 1, no week marked done yet
$ws.Range("G23").Value = 1

# Row 24: total = 1, no week marked done yet
$ws.Range("G24").Value = 1

# Move selection to reflect where the user ended up editing
$ws.Range("G26").Select()

# Reorder the header merged cells (AZ4:BA4, AO4:AP4, AR4:AS4, AU4:AV4, AX4:AY4 first)
$allMerges = @("H4:I4","K4:L4","N4:O4","Q4:R4","T4:U4","W4:X4","Z4:AA4","AC4:AD4","AF4:AG4","AI4:AJ4","AL4:AM4","AO4:AP4","AR4:AS4","AU4:AV4","AX4:AY4","AZ4:BA4")
foreach ($r in $allMerges) {
    $ws.Range($r).UnMerge()
}
$newMergeOrder = @("AZ4:BA4","AO4:AP4","AR4:AS4","AU4:AV4","AX4:AY4","AL4:AM4","H4:I4","K4:L4","N4:O4","Q4:R4","T4:U4","W4:X4","Z4:AA4","AC4:AD4","AF4:AG4","AI4:AJ4")
foreach ($r in $newMergeOrder) {
    $ws.Range($r).Merge()
}
